# Clean up non-public new spells
# Adds 9 new playtest spells to the Spells sheet, and switches the
# active/selected sheet from Spells to Subclasses.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spells")

$rows = @(
    @("Binding Chain",    "New", 1, "Conjuration",  "Yes","No","No","No","Yes","Yes","No","Yes","No"),
    @("Molten Sphere",    "New", 6, "Conjuration",  "No","No","No","Yes","No","No","Yes","No","No"),
    @("Mud Ball",         "New", 2, "Conjuration",  "No","No","No","Yes","No","Yes","No","No","No"),
    @("Tranquility",      "New", 5, "Conjuration",  "No","No","No","Yes","No","No","No","No","No"),
    @("Frozen Tomb",      "New", 5, "Evocation",    "No","No","No","Yes","No","No","No","Yes","No"),
    @("Water Whip",       "New", 1, "Conjuration",  "No","No","No","Yes","No","No","No","No","No"),
    @("Grasping Tide",    "New", 2, "Conjuration",  "No","No","No","Yes","No","No","No","No","No"),
    @("Stream of Flames", "New", 6, "Evocation",    "No","No","No","No","No","No","Yes","No","Yes"),
    @("Leap Slam",        "New", 3, "Trasmutation", "No","No","No","No","Yes","No","No","No","No")
)

$startRow = 39
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
    $ws.Cells.Item($r, 14).Value = "Playtest Ready"
    $ws.Cells.Item($r, 15).Value = "Not Released"
}

# Switch the selected/active sheet from Spells to Subclasses.
$subclasses = $wb.Worksheets.Item("Subclasses")
$subclasses.Activate()
